$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 0.6658792748069419
$ws.Range("C2").Value = 0.665879274806943
$ws.Range("D2").Value = 0.665879274806943

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.03542861843485547
$ws.Range("C3").Value = 0.03195324278056939
$ws.Range("D3").Value = 0.1320174499681472

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03587529838346575
$ws.Range("C4").Value = 0.03389033849072908
$ws.Range("D4").Value = 0.05238001152482729

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.03836556192561748
$ws.Range("C5").Value = 0.05306162557692551
$ws.Range("D5").Value = 0.06921077816248632
